$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1961.5454
$ws.Range("I43").Value = 895.25
$ws.Range("J43").Value = 2570.8572
$ws.Range("K43").Value = 895.25
$ws.Range("L43").Value = 2570.8572
$ws.Range("M43").Value = -826.25
$ws.Range("N43").Value = -2708.8572
$ws.Range("H81").Value = 30200
$ws.Range("J81").Value = 30200
$ws.Range("L81").Value = 30200
$ws.Range("N81").Value = -32196
$ws.Range("H84").Value = 30200
$ws.Range("J84").Value = 30200
$ws.Range("L84").Value = 90600
$ws.Range("N84").Value = -100584
$ws.Range("H135").Value = 669.8461
$ws.Range("I135").Value = 292.5
$ws.Range("K135").Value = 2632.5
$ws.Range("M135").Value = -97.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1287.5714
$ws.Range("I2").Value = 1400
$ws.Range("J2").Value = 1137.6666
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1137.6666
$ws.Range("M2").Value = -1287
$ws.Range("N2").Value = -1363.6666
$ws.Range("H32").Value = 5251.508
$ws.Range("I32").Value = 3821.0942
$ws.Range("K32").Value = 3821.0942
$ws.Range("M32").Value = -3534.0942
$ws.Range("H61").Value = 1102.2413
$ws.Range("I61").Value = 936.3889
$ws.Range("J61").Value = 1373.6364
$ws.Range("K61").Value = 936.3889
$ws.Range("L61").Value = 1373.6364
$ws.Range("M61").Value = -724.3889
$ws.Range("N61").Value = -1797.6364
$ws.Range("H74").Value = 2260.93
$ws.Range("I74").Value = 2199.261
$ws.Range("K74").Value = 2199.261
$ws.Range("M74").Value = -1325.261
$ws.Range("H77").Value = 2260.93
$ws.Range("I77").Value = 2199.261
$ws.Range("K77").Value = 10996.305
$ws.Range("M77").Value = -6628.305
$ws.Range("H80").Value = 25353.834
$ws.Range("J80").Value = 25353.834
$ws.Range("L80").Value = 25353.834
$ws.Range("N80").Value = -27349.834
$ws.Range("H83").Value = 25353.834
$ws.Range("J83").Value = 25353.834
$ws.Range("L83").Value = 76061.50199999999
$ws.Range("N83").Value = -86045.50199999999
$ws.Range("H116").Value = 1287.5714
$ws.Range("I116").Value = 1400
$ws.Range("J116").Value = 1137.6666
$ws.Range("K116").Value = 1400
$ws.Range("L116").Value = 1137.6666
$ws.Range("M116").Value = 894
$ws.Range("N116").Value = -5725.6666
$ws.Range("H136").Value = 1102.2413
$ws.Range("I136").Value = 936.3889
$ws.Range("J136").Value = 1373.6364
$ws.Range("K136").Value = 2809.1667
$ws.Range("L136").Value = 4120.9092
$ws.Range("M136").Value = -259.1667000000002
$ws.Range("N136").Value = -9220.9092
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1287.5714
$ws.Range("I3").Value = 1400
$ws.Range("J3").Value = 1137.6666
$ws.Range("K3").Value = 1400
$ws.Range("L3").Value = 1137.6666
$ws.Range("M3").Value = -1286
$ws.Range("N3").Value = -1365.6666
$ws.Range("H105").Value = 1811.5
$ws.Range("I105").Value = 1666.55
$ws.Range("K105").Value = 1666.55
$ws.Range("M105").Value = 80.45000000000005
$ws.Range("H107").Value = 1755
$ws.Range("I107").Value = 1719.5
$ws.Range("J107").Value = 1790.5
$ws.Range("K107").Value = 1719.5
$ws.Range("L107").Value = 1790.5
$ws.Range("M107").Value = 200.5
$ws.Range("N107").Value = -5630.5
$ws.Range("H134").Value = 2532.2273
$ws.Range("I134").Value = 1414.7693
$ws.Range("J134").Value = 4146.3335
$ws.Range("K134").Value = 4244.3079
$ws.Range("L134").Value = 12439.0005
$ws.Range("M134").Value = -1709.3079
$ws.Range("N134").Value = -17509.0005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8335403
$ws.Range("I31").Value = 1154.2559
$ws.Range("J31").Value = 29416148
$ws.Range("K31").Value = 1154.2559
$ws.Range("L31").Value = 29416148
$ws.Range("M31").Value = -859.2559000000001
$ws.Range("N31").Value = -29416738
$ws.Range("H34").Value = 8335403
$ws.Range("I34").Value = 1154.2559
$ws.Range("J34").Value = 29416148
$ws.Range("K34").Value = 1154.2559
$ws.Range("L34").Value = 29416148
$ws.Range("M34").Value = -952.2559000000001
$ws.Range("N34").Value = -29416552
$ws.Range("H134").Value = 3242.4182
$ws.Range("I134").Value = 3467.1892
$ws.Range("J134").Value = 2780.389
$ws.Range("K134").Value = 10401.5676
$ws.Range("L134").Value = 8341.167000000001
$ws.Range("M134").Value = -7866.567599999998
$ws.Range("N134").Value = -13411.167
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1402.919
$ws.Range("I5").Value = 487.1875
$ws.Range("J5").Value = 2100.6191
$ws.Range("K5").Value = 1461.5625
$ws.Range("L5").Value = 6301.8573
$ws.Range("M5").Value = -1349.5625
$ws.Range("N5").Value = -6525.8573
$ws.Range("H92").Value = 402.42856
$ws.Range("J92").Value = 453
$ws.Range("L92").Value = 1359
$ws.Range("N92").Value = -3855
$ws.Range("H96").Value = 7000
$ws.Range("J96").Value = 7000
$ws.Range("L96").Value = 21000
$ws.Range("N96").Value = -25118
$ws.Range("H107").Value = 31728.25
$ws.Range("I107").Value = 544.5833
$ws.Range("J107").Value = 50438.45
$ws.Range("K107").Value = 1633.7499
$ws.Range("L107").Value = 151315.35
$ws.Range("M107").Value = 286.2501
$ws.Range("N107").Value = -155155.35
$ws.Range("H135").Value = 1402.919
$ws.Range("I135").Value = 487.1875
$ws.Range("J135").Value = 2100.6191
$ws.Range("K135").Value = 4384.6875
$ws.Range("L135").Value = 18905.5719
$ws.Range("M135").Value = -1849.6875
$ws.Range("N135").Value = -23975.5719
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 22500
$ws.Range("J27").Value = 22500
$ws.Range("L27").Value = 22500
$ws.Range("N27").Value = -22832
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H80").Value = 15627344
$ws.Range("I80").Value = 62501400
$ws.Range("J80").Value = 2658.3333
$ws.Range("K80").Value = 62501400
$ws.Range("L80").Value = 2658.3333
$ws.Range("M80").Value = -62500402
$ws.Range("N80").Value = -4654.3333
$ws.Range("H83").Value = 15627344
$ws.Range("I83").Value = 62501400
$ws.Range("J83").Value = 2658.3333
$ws.Range("K83").Value = 312507000
$ws.Range("L83").Value = 13291.6665
$ws.Range("M83").Value = -312502008
$ws.Range("N83").Value = -23275.6665
$ws.Range("H102").Value = 2303.8125
$ws.Range("I102").Value = 1427.3636
$ws.Range("J102").Value = 4232
$ws.Range("K102").Value = 1427.3636
$ws.Range("L102").Value = 4232
$ws.Range("M102").Value = 194.6364000000001
$ws.Range("N102").Value = -7476
$ws.Range("H122").Value = 5842.857
$ws.Range("I122").Value = 2633.3333
$ws.Range("J122").Value = 8250
$ws.Range("K122").Value = 7899.999899999999
$ws.Range("L122").Value = 24750
$ws.Range("M122").Value = -5449.999899999999
$ws.Range("N122").Value = -29650
$ws.Range("H132").Value = 1905.7451
$ws.Range("I132").Value = 1109.8158
$ws.Range("J132").Value = 4232.3076
$ws.Range("K132").Value = 3329.4474
$ws.Range("L132").Value = 12696.9228
$ws.Range("M132").Value = -799.4474
$ws.Range("N132").Value = -17756.9228
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6020.4
$ws.Range("I40").Value = 4449.909
$ws.Range("K40").Value = 4449.909
$ws.Range("M40").Value = -4313.909
$ws.Range("H122").Value = 6298
$ws.Range("I122").Value = 2758
$ws.Range("K122").Value = 8274
$ws.Range("M122").Value = -5824
$ws.Range("H132").Value = 8048.14
$ws.Range("I132").Value = 7432.9736
$ws.Range("J132").Value = 9996.166999999999
$ws.Range("K132").Value = 22298.9208
$ws.Range("L132").Value = 29988.501
$ws.Range("M132").Value = -19768.9208
$ws.Range("N132").Value = -35048.501
$ws.Range("H136").Value = 1850.2222
$ws.Range("I136").Value = 880.4888999999999
$ws.Range("J136").Value = 6698.8887
$ws.Range("K136").Value = 2641.4667
$ws.Range("L136").Value = 20096.6661
$ws.Range("M136").Value = -91.46669999999995
$ws.Range("N136").Value = -25196.6661
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 6456
$ws.Range("I26").Value = 2912
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 2912
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = -2619
$ws.Range("N26").Value = -10586
$ws.Range("H29").Value = 34173.668
$ws.Range("I29").Value = 16255
$ws.Range("J29").Value = 70011
$ws.Range("K29").Value = 16255
$ws.Range("L29").Value = 70011
$ws.Range("M29").Value = -15965
$ws.Range("N29").Value = -70591
$ws.Range("H136").Value = 2352.3096
$ws.Range("I136").Value = 720.70966
$ws.Range("K136").Value = 2162.12898
$ws.Range("M136").Value = 387.87102
